$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exportar")

# Update shared date strings in column G (Data de Inicio) and H (Data Final)
# These strings are shared across all data rows (2-19)
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 7).Value = "14/07/2025 06:00"
    $ws.Cells.Item($r, 8).Value = "15/07/2025 06:00"
}

# Update numeric metric values per row
# Row 2
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 8640.9
$ws.Range("K2").Value = 47.172780555555555
$ws.Range("L2").Value = 143.31858222222223
$ws.Range("M2").Value = 0.31454305555555556
$ws.Range("N2").Value = 6.329743471888155
$ws.Range("O2").Value = 20.096036440258775
$ws.Range("P2").Value = 35.266
$ws.Range("Q2").Value = 190.8059072222222
$ws.Range("R2").Value = 629.5660922345413
$ws.Range("S2").Value = 0.34
$ws.Range("T2").Value = 13.103523297333458
$ws.Range("U2").Value = 43.23521256183475
$ws.Range("V2").Value = 0.26527062450343325

# Row 4
$ws.Range("J4").Value = 8254.45

# Row 7
$ws.Range("I7").Value = 10.649999999999636
$ws.Range("J7").Value = 8891
$ws.Range("K7").Value = 22.32755166666667
$ws.Range("L7").Value = 96.25951916666666
$ws.Range("M7").Value = 1.7534527777777777
$ws.Range("N7").Value = 5.7279570056632725
$ws.Range("O7").Value = 14.394394104136987
$ws.Range("P7").Value = 20.347662779397474
$ws.Range("Q7").Value = 120.34052472222223
$ws.Range("R7").Value = 397.064823567475
$ws.Range("S7").Value = 0.64
$ws.Range("T7").Value = 11.317118837813531
$ws.Range("U7").Value = 37.340952310130476
$ws.Range("V7").Value = 0.2754658760152494

# Row 8
$ws.Range("I8").Value = 20.25
$ws.Range("J8").Value = 6938.6
$ws.Range("K8").Value = 48.28416805555555
$ws.Range("L8").Value = 184.4372263888889
$ws.Range("M8").Value = 15.846829722222223
$ws.Range("N8").Value = 4.697706192153328
$ws.Range("O8").Value = 18.648808259421326
$ws.Range("P8").Value = 31.132404060248856
$ws.Range("Q8").Value = 248.56822444444447
$ws.Range("R8").Value = 820.1534637756838
$ws.Range("S8").Value = 0.16
$ws.Range("T8").Value = 12.058364336364917
$ws.Range("U8").Value = 39.7866996074924
$ws.Range("V8").Value = 0.26047755264962713

# Row 9
$ws.Range("I9").Value = 19.600000000000364
$ws.Range("J9").Value = 7751.15
$ws.Range("K9").Value = 45.92582499999999
$ws.Range("L9").Value = 197.8471436111111
$ws.Range("M9").Value = 16.195511944444444
$ws.Range("N9").Value = 4.681169434844015
$ws.Range("O9").Value = 19.77231551493919
$ws.Range("P9").Value = 32.579031690140845
$ws.Range("Q9").Value = 259.96848388888884
$ws.Range("R9").Value = 857.7687393894508
$ws.Range("S9").Value = 0.18
$ws.Range("T9").Value = 12.851220712666365
$ws.Range("U9").Value = 42.40273753733489
$ws.Range("V9").Value = 0.2601261926095682

# Row 10
$ws.Range("I10").Value = 14.600000000000364
$ws.Range("J10").Value = 9340.15
$ws.Range("K10").Value = 21.80077333333333
$ws.Range("L10").Value = 134.8915322222222
$ws.Range("M10").Value = 6.5919847222222225
$ws.Range("N10").Value = 3.9907509080481747
$ws.Range("O10").Value = 15.014758835585312
$ws.Range("P10").Value = 29.57345411352838
$ws.Range("Q10").Value = 163.2842894444444
$ws.Range("R10").Value = 538.7582257037178
$ws.Range("S10").Value = 0.5479999999999999
$ws.Range("T10").Value = 11.172662548234607
$ws.Range("U10").Value = 36.864317267469104
$ws.Range("V10").Value = 0.2647239987289581

# Row 11
$ws.Range("I11").Value = 13.149999999999636
$ws.Range("J11").Value = 8720.1
$ws.Range("K11").Value = 42.97993277777777
$ws.Range("L11").Value = 106.2052063888889
$ws.Range("M11").Value = 1.641596111111111
$ws.Range("N11").Value = 5.7720673693152555
$ws.Range("O11").Value = 18.51605320040093
$ws.Range("P11").Value = 33.359700170357755
$ws.Range("Q11").Value = 150.82673833333334
$ws.Range("R11").Value = 497.65440514589784
$ws.Range("S11").Value = 0.508
$ws.Range("T11").Value = 11.423222839710496
$ws.Range("U11").Value = 37.69104357730974
$ws.Range("V11").Value = 0.2860285416571131

# Row 12
$ws.Range("I12").Value = 18.75
$ws.Range("J12").Value = 1297.7
$ws.Range("K12").Value = 93.36399777777775
$ws.Range("L12").Value = 236.47190638888893
$ws.Range("M12").Value = 0.7499136111111111
$ws.Range("N12").Value = 10.248092057046346
$ws.Range("O12").Value = 19.388185992703487
$ws.Range("P12").Value = 25.13116713352007
$ws.Range("Q12").Value = 330.58581527777784
$ws.Range("R12").Value = 1090.7713650091919
$ws.Range("S12").Value = 0.20800000000000002
$ws.Range("T12").Value = 15.493474170104363
$ws.Range("U12").Value = 51.12088053463257
$ws.Range("V12").Value = 0.3459918551635875

# Row 13
$ws.Range("I13").Value = 17.049999999999955
$ws.Range("J13").Value = 1285
$ws.Range("K13").Value = 49.020275555555564
$ws.Range("L13").Value = 171.1249558333333
$ws.Range("M13").Value = 1.5798097222222223
$ws.Range("N13").Value = 6.616403682517787
$ws.Range("O13").Value = 16.87365295172886
$ws.Range("P13").Value = 29.913199472295513
$ws.Range("Q13").Value = 221.72504277777784
$ws.Range("R13").Value = 731.5841043095583
$ws.Range("S13").Value = 0.428
$ws.Range("T13").Value = 12.595245724951639
$ws.Range("U13").Value = 41.5581453804597
$ws.Range("V13").Value = 0.30565065841301353

# Row 19
$ws.Range("I19").Value = 18.350000000000136
$ws.Range("J19").Value = 1383.45
$ws.Range("K19").Value = 94.83102361111112
$ws.Range("L19").Value = 183.82474277777774
$ws.Range("M19").Value = 1.0984266666666667
$ws.Range("N19").Value = 10.753571261645444
$ws.Range("O19").Value = 19.354920152891875
$ws.Range("P19").Value = 24.31918411330049
$ws.Range("Q19").Value = 279.754191111111
$ws.Range("R19").Value = 923.051887900591
$ws.Range("S19").Value = 0.364
$ws.Range("T19").Value = 15.23613638928345
$ws.Range("U19").Value = 50.271791827609235
$ws.Range("V19").Value = 0.3941620680058883
